# Weekly data refresh: insert a new week's reading as row 10 (pushing the
# existing rows 10-15 down to 11-16) and populate it with the latest
# Espárragos / Vega Modelo de Temuco reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 10..15 down to 11..16, growing the used range to A1:R16.
$ws.Rows.Item(10).Insert()

# Fill the newly-opened row 10 with the new week's record.
$ws.Cells.Item(10, 1).Value = 10
$ws.Cells.Item(10, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(10, 3).Value = "La Araucanía"
$ws.Cells.Item(10, 4).Value = 44463
$ws.Cells.Item(10, 5).Value = 9
$ws.Cells.Item(10, 6).Value = 300000000
$ws.Cells.Item(10, 7).Value = "Espárragos"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 40
$ws.Cells.Item(10, 11).Value = 2500
$ws.Cells.Item(10, 12).Value = 2500
$ws.Cells.Item(10, 13).Value = 2500
$ws.Cells.Item(10, 14).Value = "$/kilo"
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 2500
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = "Hortaliza"
